$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '55.311.31'
$ws.Range("E2").Value = '  -2.61%  '
$ws.Range("D3").Value = '2.914.31'
$ws.Range("E3").Value = '  -1.73%  '
$ws.Range("E4").Value = '  +0.67%  '
$ws.Range("D5").Value2 = "'495.08"
$ws.Range("E5").Value = '  -0.35%  '
$ws.Range("D6").Value2 = "'134.63"
$ws.Range("E6").Value = '  -1.49%  '
$ws.Range("D7").Value2 = "'1.00"
$ws.Range("E7").Value = '  +0.41%  '
$ws.Range("D8").Value2 = "'0.423"
$ws.Range("E8").Value = '  -0.23%  '
$ws.Range("D9").Value2 = "'6.97"
$ws.Range("E9").Value = '  -4.40%  '
$ws.Range("E10").Value = '  -1.41%  '
$ws.Range("D11").Value2 = "'0.360"
$ws.Range("E11").Value = '  +1.41%  '
$ws.Range("D12").Value = '3.515.53'
$ws.Range("E12").Value = '  +1.16%  '
$ws.Range("E13").Value = '  -2.64%  '
$ws.Range("D14").Value2 = "'25.66"
$ws.Range("E14").Value = '  -0.11%  '
$ws.Range("D15").Value2 = "'0.0000158"
$ws.Range("E15").Value = '  +0.99%  '
$ws.Range("D16").Value = '55.697.22'
$ws.Range("E16").Value = '  -2.09%  '
$ws.Range("D17").Value = '2.998.27'
$ws.Range("E17").Value = '  +1.27%  '
$ws.Range("D18").Value2 = "'5.85"
$ws.Range("E18").Value = '  -2.89%  '
$ws.Range("D19").Value2 = "'12.74"
$ws.Range("E19").Value = '  +1.45%  '
$ws.Range("D20").Value2 = "'7.71"
$ws.Range("E20").Value = '  -0.86%  '
$ws.Range("D21").Value2 = "'319.71"
$ws.Range("E21").Value = '  +0.74%  '
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("D23").Value2 = "'0.483"
$ws.Range("E23").Value = '  -0.14%  '
$ws.Range("D24").Value2 = "'63.61"
$ws.Range("E24").Value = '  +1.16%  '
$ws.Range("D25").Value = '3.144.28'
$ws.Range("E25").Value = '  +1.90%  '
$ws.Range("E26").Value = '  +1.57%  '
$ws.Range("D27").Value2 = "'0.159"
$ws.Range("E27").Value = '  -1.65%  '
$ws.Range("D28").Value = '0.0₃0857'
$ws.Range("E28").Value = '  -2.96%  '
$ws.Range("D29").Value2 = "'6.35"
$ws.Range("E29").Value = '  -2.22%  '
$ws.Range("D30").Value2 = "'6.87"
$ws.Range("E30").Value = '  -2.45%  '
$ws.Range("D31").Value2 = "'1.76"
$ws.Range("E31").Value = '  -0.23%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").Value2 = "'19.83"
$ws.Range("E32").Value = '  -1.22%  '
$ws.Range("B33").Value = 'Fetch.AI'
$ws.Range("C33").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D33").Value2 = "'1.13"
$ws.Range("E33").Value = '  -1.23%  '
$ws.Range("D34").Value2 = "'150.00"
$ws.Range("E34").Value = '  -2.81%  '
$ws.Range("D35").Value2 = "'4.48"
$ws.Range("E35").Value = '  -2.30%  '
$ws.Range("D36").Value2 = "'5.66"
$ws.Range("E36").Value = '  -0.55%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value2 = "'1.22"
$ws.Range("E37").Value = '  -1.59%  '
$ws.Range("B38").Value = 'EnergySwap'
$ws.Range("C38").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D38").Value2 = "'24.44"
$ws.Range("E38").Value = '  +2.41%  '
$ws.Range("D39").Value2 = "'0.0647"
$ws.Range("E39").Value = '  -2.28%  '
$ws.Range("D40").Value = '3.014.71'
$ws.Range("E40").Value = '  +0.56%  '
$ws.Range("E41").Value = '  +0.71%  '
$ws.Range("D42").Value2 = "'36.51"
$ws.Range("E42").Value = '  -2.45%  '
$ws.Range("D43").Value2 = "'0.649"
$ws.Range("E43").Value = '  +2.00%  '
$ws.Range("D44").Value2 = "'3.68"
$ws.Range("E44").Value = '  -0.39%  '
$ws.Range("D45").Value = '2.139.42'
$ws.Range("E45").Value = '  -2.41%  '
$ws.Range("E46").Value = '  -3.60%  '
$ws.Range("B47").Value = 'Cosmos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D47").Value2 = "'5.79"
$ws.Range("E47").Value = '  -2.00%  '
$ws.Range("B48").Value = 'ONDO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D48").Value2 = "'0.913"
$ws.Range("E48").Value = '  -2.36%  '
$ws.Range("D49").Value2 = "'0.0234"
$ws.Range("E49").Value = '  +0.18%  '
$ws.Range("D50").Value2 = "'19.31"
$ws.Range("E50").Value = '  +1.40%  '
$ws.Range("D51").Value2 = "'0.0838"
$ws.Range("E51").Value = '  -4.20%  '
